{"js": "// Site rebuild: drop the trailing \"Ver no Jupiter / Salvar em pdf / Salvar\n// em docx\" line, the \"\u00a9 2020 ...\" footer line, and the blank paragraph that\n// separated them from the Bibliografia text.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter...\" paragraph; the blank paragraph right before\n// it and the copyright paragraph right after it are removed together with\n// it.\nlet verIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex !== -1) {\n  const toDelete = [];\n  const prev = items[verIndex - 1];\n  if (prev && prev.text.trim() === \"\") {\n    toDelete.push(prev);\n  }\n  toDelete.push(items[verIndex]);\n  const next = items[verIndex + 1];\n  if (\n    next &&\n    next.text.trim() ===\n      \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n  ) {\n    toDelete.push(next);\n  }\n  for (const p of toDelete) {\n    p.delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Site rebuild: drop the trailing \"Ver no Jupiter / Salvar em pdf / Salvar\n# em docx\" line, the \"(c) 2020 ...\" footer line, and the blank paragraph\n# that separated them from the Bibliografia text.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0xA9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$verIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -ge 1) {\n    $nextIndex = $verIndex + 1\n    if ($nextIndex -le $count) {\n        $nextText = $d.Paragraphs.Item($nextIndex).Range.Text.Trim()\n        if ($nextText -eq $copyrightText) {\n            $d.Paragraphs.Item($nextIndex).Range.Delete()\n        }\n    }\n\n    $d.Paragraphs.Item($verIndex).Range.Delete()\n\n    $prevIndex = $verIndex - 1\n    if ($prevIndex -ge 1) {\n        $prevText = $d.Paragraphs.Item($prevIndex).Range.Text.Trim()\n        if ($prevText -eq \"\") {\n            $d.Paragraphs.Item($prevIndex).Range.Delete()\n        }\n    }\n}\n"}
